# Adds a new "2022-Q1" sheet (holdings detail) before the existing "总计"
# (totals) sheet, and inserts a corresponding summary row at the top of the
# "总计" sheet's data table.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet immediately before "总计".
# ------------------------------------------------------------------
$totalSheet = $sheets.Item($sheets.Count)
$newSheet = $sheets.Add($totalSheet, $null)
$newSheet.Name = "2022-Q1"

# Header row (B1:H1)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Data rows: index|code|name|scale|position|ratio|marketvalue|rank
$data = @"
0|512880|国泰中证全指证券公司ETF|322.34|99.85|2.73|8.7999|10
1|512000|华宝中证全指证券公司ETF|230.47|99.86|2.71|6.2457|10
2|512900|南方中证全指证券公司ETF|78.29|99.89|2.71|2.1217|10
3|159841|天弘中证全指证券公司ETF|42.72|99.81|2.71|1.1577|10
4|161720|招商中证全指证券公司指数（LOF）|28.99|94.46|2.56|0.7421|10
5|004814|中欧红利优享灵活配置混合A|22.96|93.96|2.73|0.6268|9
6|501016|国泰中证申万证券行业指数（LOF）|19.81|93.96|2.60|0.5151|10
7|501047|汇添富中证全指证券公司指数LOF A|18.61|94.26|2.62|0.4876|10
8|163113|申万菱信中证申万证券行业指数（LOF）|18.63|93.17|2.56|0.4769|10
9|160633|鹏华中证全指证券公司指数（LOF）|17.69|94.30|2.55|0.4511|10
10|159993|鹏华国证证券龙头ETF|13.43|97.76|3.08|0.4136|10
11|161027|富国中证全指证券公司指数|14.89|94.31|2.58|0.3842|10
12|502010|易方达证券公司指数（LOF）|14.03|94.56|2.56|0.3592|10
13|515010|华夏中证全指证券公司ETF|11.46|99.36|2.69|0.3083|10
14|501048|汇添富中证全指证券公司指数LOF C|10.39|94.26|2.62|0.2722|10
15|004496|前海开源多元策略灵活配置混合A|3.09|91.25|7.64|0.2361|6
16|004815|中欧红利优享灵活配置混合C|8.45|93.96|2.73|0.2307|9
17|160516|博时中证全指证券公司指数|8.30|94.96|2.56|0.2125|10
18|004497|前海开源多元策略灵活配置混合C|2.06|91.25|7.64|0.1574|6
19|005541|前海开源盛鑫灵活配置混合A|2.74|94.03|5.36|0.1469|10
20|160419|华安中证全指证券公司指数（LOF）A|5.41|94.03|2.57|0.1390|10
21|515560|建信中证全指证券公司ETF|4.71|98.22|2.72|0.1281|10
22|159842|银华中证全指证券公司ETF|4.73|97.53|2.68|0.1268|10
23|502053|长盛中证全指证券公司指数（LOF）|4.28|94.45|2.55|0.1091|10
24|515850|富国中证全指证券公司ETF|2.17|99.49|2.70|0.0586|10
25|159848|国联安中证全指证券公司ETF|2.13|98.30|2.68|0.0571|10
26|011997|景顺长城安盈回报一年持有期混合型证券投资基金A|5.74|29.73|0.96|0.0551|10
27|512570|易方达中证全指证券公司ETF|1.74|99.11|2.69|0.0468|10
28|011243|万家惠裕回报6个月持有期混合型证券投资基金A|4.93|23.04|0.69|0.0340|10
29|000423|前海开源事件驱动混合A|0.48|88.63|6.77|0.0325|7
30|510200|汇安上证证券ETF|0.74|97.70|4.34|0.0321|7
31|008116|银华沪深股通精选混合|0.55|88.86|3.64|0.0200|4
32|005444|光大保德信多策略精选18个月定期开放灵活配置混合|1.00|29.09|1.97|0.0197|7
33|000417|国联安新精选灵活配置混合|0.52|28.16|2.13|0.0111|3
34|005542|前海开源盛鑫灵活配置混合C|0.17|94.03|5.36|0.0091|10
35|516200|华安中证全指证券公司ETF|0.25|93.02|2.54|0.0064|10
36|001865|前海开源事件驱动混合C|0.05|88.63|6.77|0.0034|7
37|011998|景顺长城安盈回报一年持有期混合型证券投资基金C|0.23|29.73|0.96|0.0022|10
38|162907|泰信中证锐联基本面400指数（LOF）|0.23|94.61|0.54|0.0012|8
39|011244|万家惠裕回报6个月持有期混合型证券投资基金C|0.14|23.04|0.69|0.0010|10
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $idx = [int]$parts[0]
    $code = $parts[1]
    $name = $parts[2]
    $scale = $parts[3]
    $position = $parts[4]
    $ratio = $parts[5]
    $marketValue = $parts[6]
    $rank = [int]$parts[7]

    $row = $idx + 2

    # Column A: numeric row index, bold + bordered (matches other sheets)
    $cellA = $newSheet.Cells.Item($row, 1)
    $cellA.Value = $idx
    $cellA.Font.Bold = $true
    $cellA.Borders.LineStyle = 1
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160

    # Columns B-G: text values (force text so numeric-looking strings and
    # leading-zero codes are preserved verbatim, matching the source data).
    $cellB = $newSheet.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $code

    $cellC = $newSheet.Cells.Item($row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $name

    $cellD = $newSheet.Cells.Item($row, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $scale

    $cellE = $newSheet.Cells.Item($row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $position

    $cellF = $newSheet.Cells.Item($row, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $ratio

    $cellG = $newSheet.Cells.Item($row, 7)
    $cellG.NumberFormat = "@"
    $cellG.Value = $marketValue

    # Column H: numeric rank
    $newSheet.Cells.Item($row, 8).Value = $rank
}

# ------------------------------------------------------------------
# 2. Insert a new "2022-Q1" summary row at the top of the "总计" sheet's
#    data table (row 2), pushing the other quarters down by one row.
#
# NOTE: $totalSheet was captured *before* the new sheet was inserted in
# front of it, so its cached .Index is now stale (still reports its old
# position). Re-resolve the worksheet by name so we operate on the
# actual "总计" sheet at its new position.
# ------------------------------------------------------------------
$totalSheet = $sheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The A column is a running 0-based row counter, so every pre-existing
# row's counter value also shifts up by one (old row 2 had A=0, it is now
# row 3 and must read A=1, and so on through the old final row).
for ($r = 7; $r -ge 3; $r--) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

$cellA2 = $totalSheet.Cells.Item(2, 1)
$cellA2.Value = 0
$cellA2.Font.Bold = $true
$cellA2.Borders.LineStyle = 1
$cellA2.HorizontalAlignment = -4108
$cellA2.VerticalAlignment = -4160

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 40
$totalSheet.Cells.Item(2, 4).Value = 25.24
